$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$sheet1Changes = @{
  4  = 9339
  5  = 579
  8  = 256
  9  = 342
  10 = 400
  14 = 429
  15 = 11951
  17 = 306
  22 = 170
  23 = 115
  24 = 155
  25 = 2716
  26 = 2095
  30 = 2142
  31 = 988
  32 = 4184
  33 = 3611
  34 = 470
  37 = 13
  38 = 1310
  39 = 190
  40 = 771
  42 = 410
  43 = 493
  46 = 212
  47 = 108
  49 = 130
}
foreach ($row in $sheet1Changes.Keys) {
  $ws1.Cells.Item($row, 6).Value = $sheet1Changes[$row]
}

# --- Sheet 2: 演出 ---
$ws2 = $wb.Worksheets.Item("演出")
$sheet2Changes = @{
  2  = 20
  14 = 23
}
foreach ($row in $sheet2Changes.Keys) {
  $ws2.Cells.Item($row, 6).Value = $sheet2Changes[$row]
}

# --- Sheet 4: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$sheet4Changes = @{
  2  = 20
  7  = 9339
  8  = 579
  11 = 256
  12 = 342
  13 = 400
  15 = 429
  16 = 11951
  22 = 170
  23 = 115
  24 = 155
  25 = 2716
  26 = 2095
  31 = 2142
  32 = 988
  33 = 4184
  34 = 3611
  35 = 470
  38 = 13
  39 = 1310
  40 = 190
  41 = 771
  42 = 410
  43 = 493
  46 = 212
  47 = 108
  49 = 130
}
foreach ($row in $sheet4Changes.Keys) {
  $ws4.Cells.Item($row, 6).Value = $sheet4Changes[$row]
}
